$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$CellRef,
        [string]$TextValue
    )
    $c = $ws.Range($CellRef)
    $c.NumberFormat = "@"
    $c.Value = $TextValue
    $c.Style = "Normal"
}

Set-CellText "D2" '63.615.61'
$ws.Range("E2").Value = '  -1.26%  '
Set-CellText "D3" '3.083.81'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  +0.05%  '
Set-CellText "D5" '554.00'
$ws.Range("E5").Value = '  +0.26%  '
Set-CellText "D6" '136.33'
$ws.Range("E6").Value = '  -3.85%  '
Set-CellText "D7" '1.00'
$ws.Range("E7").Value = '  +0.04%  '
Set-CellText "D8" '3.073.54'
$ws.Range("E8").Value = '  +0.07%  '
Set-CellText "D9" '0.490'
$ws.Range("E9").Value = '  -0.25%  '
Set-CellText "D10" '6.59'
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("E11").Value = '  +1.77%  '
Set-CellText "D12" '0.447'
$ws.Range("E12").Value = '  -2.37%  '
Set-CellText "D13" '34.78'
$ws.Range("E13").Value = '  -3.76%  '
Set-CellText "D14" '0.0000214'
$ws.Range("E14").Value = '  -1.87%  '
Set-CellText "D15" '3.581.35'
$ws.Range("E15").Value = '  +0.22%  '
Set-CellText "D16" '63.733.09'
$ws.Range("E16").Value = '  -1.09%  '
$ws.Range("E17").Value = '  -0.41%  '
Set-CellText "D18" '3.083.50'
$ws.Range("E18").Value = '  -0.02%  '
Set-CellText "D19" '502.33'
$ws.Range("E19").Value = '  +2.29%  '
Set-CellText "D20" '6.56'
$ws.Range("E20").Value = '  -1.52%  '
Set-CellText "D21" '13.46'
$ws.Range("E21").Value = '  -1.76%  '
Set-CellText "D22" '0.696'
$ws.Range("E22").Value = '  +0.92%  '
Set-CellText "D23" '7.15'
$ws.Range("E23").Value = '  -1.77%  '
Set-CellText "D24" '12.25'
$ws.Range("E24").Value = '  -1.76%  '
Set-CellText "D25" '76.76'
$ws.Range("E25").Value = '  -2.37%  '
Set-CellText "D26" '0.997'
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("E27").Value = '  +1.19%  '
Set-CellText "D28" '8.19'
$ws.Range("E28").Value = '  +1.35%  '
Set-CellText "D29" '2.01'
$ws.Range("E29").Value = '  -3.75%  '
$ws.Range("E30").Value = '  -0.15%  '
Set-CellText "D31" '25.95'
$ws.Range("E31").Value = '  -0.04%  '
Set-CellText "D32" '2.50'
$ws.Range("E32").Value = '  -5.78%  '
Set-CellText "D33" '1.10'
$ws.Range("E33").Value = '  -3.60%  '
Set-CellText "D34" '528.63'
$ws.Range("E34").Value = '  -10.84%  '
Set-CellText "D35" '57.71'
$ws.Range("E35").Value = '  +10.27%  '
Set-CellText "D36" '5.80'
$ws.Range("E36").Value = '  -3.71%  '
Set-CellText "D37" '5.09'
$ws.Range("E37").Value = '  -6.33%  '
Set-CellText "D38" '0.0408'
$ws.Range("E38").Value = '  +0.98%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText "D39" '3.044.04'
$ws.Range("E39").Value = '  +2.37%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText "D40" '0.0787'
$ws.Range("E40").Value = '  -1.24%  '
$ws.Range("E41").Value = '  -1.33%  '
Set-CellText "D42" '8.01'
$ws.Range("E42").Value = '  -3.00%  '
Set-CellText "D43" '2.53'
$ws.Range("E43").Value = '  -12.06%  '
$ws.Range("E44").Value = '  +0.06%  '
Set-CellText "D45" '0.248'
$ws.Range("E45").Value = '  +0.41%  '
Set-CellText "D46" '2.04'
$ws.Range("E46").Value = '  -3.07%  '
Set-CellText "D47" '121.97'
$ws.Range("E47").Value = '  +1.47%  '
Set-CellText "D48" '23.87'
$ws.Range("E48").Value = '  -4.77%  '
$ws.Range("E49").Value = '  -2.10%  '
Set-CellText "D50" '0.0₃0488'
$ws.Range("E50").Value = '  -9.17%  '
Set-CellText "D51" '1.99'
$ws.Range("E51").Value = '  -4.28%  '
